$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.531.24"
$ws.Range("E2").Value = "  +0.28%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.951.66"
$ws.Range("E3").Value = "  +0.36%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.71"
$ws.Range("E5").Value = "  +0.50%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.612"
$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "57.69"
$ws.Range("E7").Value = "  +1.19%  "

$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.374"
$ws.Range("E9").Value = "  +3.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0787"
$ws.Range("E10").Value = "  -7.48%  "

$ws.Range("E11").Value = "  -1.10%  "

$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.830"
$ws.Range("E12").Value = "  +1.78%  "

$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.236.24"
$ws.Range("E13").Value = "  +0.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "13.74"
$ws.Range("E14").Value = "  +1.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.04"
$ws.Range("E15").Value = "  -0.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.28"
$ws.Range("E16").Value = "  +2.11%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.950.74"
$ws.Range("E17").Value = "  +0.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.425.60"
$ws.Range("E18").Value = "  +0.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.48"
$ws.Range("E19").Value = "  +0.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0845"
$ws.Range("E20").Value = "  -3.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "229.18"
$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.01"
$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("E23").Value = "  -0.28%  "

$ws.Range("E24").Value = "  +2.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.35"
$ws.Range("E25").Value = "  +3.11%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.07"
$ws.Range("E26").Value = "  -1.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.05"
$ws.Range("E27").Value = "  -1.00%  "

$ws.Range("E28").Value = "  +0.23%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.26"
$ws.Range("E29").Value = "  +0.11%  "

$ws.Range("E30").Value = "  +1.50%  "

$ws.Range("E31").Value = "  +2.22%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.67"
$ws.Range("E32").Value = "  +1.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0607"
$ws.Range("E33").Value = "  -4.65%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.36"
$ws.Range("E34").Value = "  +2.64%  "

$ws.Range("E35").Value = "  +0.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.41"
$ws.Range("E36").Value = "  +13.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.24"
$ws.Range("E37").Value = "  +4.65%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.76"
$ws.Range("E38").Value = "  -1.45%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.20"
$ws.Range("E39").Value = "  -15.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0970"
$ws.Range("E40").Value = "  +0.10%  "

$ws.Range("E41").Value = "  +1.03%  "

$ws.Range("E42").Value = "  -0.70%  "

$ws.Range("E43").Value = "  -0.60%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.363.78"
$ws.Range("E44").Value = "  +1.41%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.70"
$ws.Range("E45").Value = "  +0.27%  "

$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.02"
$ws.Range("E46").Value = "  -0.47%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.19"
$ws.Range("E47").Value = "  +0.08%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.14"
$ws.Range("E48").Value = "  +0.71%  "

$ws.Range("E49").Value = "  +0.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.125.76"
$ws.Range("E50").Value = "  +0.17%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.66"
$ws.Range("E51").Value = "  -1.39%  "
